$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" (sheet1.xml / table3 "Overview") ---
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()

$wsOverview.Cells.Item(3, 1).Value = "9a71ffcd-7f60-46ca-8938-172a2aca64aa.md"
$wsOverview.Cells.Item(3, 2).Value = "e2e\9a71ffcd-7f60-46ca-8938-172a2aca64aa.md"
$wsOverview.Cells.Item(3, 3).Value = ".md"
$wsOverview.Cells.Item(3, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 7).Value = "2016-08-13 14:48:23"
$wsOverview.Cells.Item(3, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(3, 2), "https://github.com/OpenLocalizationTestOrg/oltest/blob/24dad09822a2e81cfe6c5f481aba677f2333a927/e2e/9a71ffcd-7f60-46ca-8938-172a2aca64aa.md", "", "", "e2e\9a71ffcd-7f60-46ca-8938-172a2aca64aa.md")

# --- Sheet "zh-cn" (sheet2.xml / table1 "zh-cn") ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$rowZh = $loZh.ListRows.Add()

$wsZh.Cells.Item(3, 1).Value = "9a71ffcd-7f60-46ca-8938-172a2aca64aa.md"
$wsZh.Cells.Item(3, 2).Value = ".md"
$wsZh.Cells.Item(3, 3).Value = "Ready for handoff"
$wsZh.Cells.Item(3, 4).Value = "e2e"
$wsZh.Cells.Item(3, 5).Value = "ht"
$wsZh.Cells.Item(3, 6).Value = "'False"
$wsZh.Cells.Item(3, 7).Value = "9a71ffcd-7f60-46ca-8938-172a2aca64aa.842dd4dc5e93682265b9bcceb6d8905f68e54fe8.zh-cn.xlf"
$wsZh.Cells.Item(3, 8).Value = "2016-08-13 14:48:15"
$wsZh.Cells.Item(3, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(3, 11).Value = "0001-01-01 00:00:00"
$wsZh.Cells.Item(3, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(3, 13).Value = "'True"
$wsZh.Cells.Item(3, 15).Value = "'False"

$wsZh.Hyperlinks.Add($wsZh.Cells.Item(3, 1), "https://github.com/OpenLocalizationTestOrg/oltest/blob/24dad09822a2e81cfe6c5f481aba677f2333a927/e2e/9a71ffcd-7f60-46ca-8938-172a2aca64aa.md", "", "", "9a71ffcd-7f60-46ca-8938-172a2aca64aa.md")

# --- Sheet "de-de" (sheet3.xml / table2 "de-de") ---
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$rowDe = $loDe.ListRows.Add()

$wsDe.Cells.Item(3, 1).Value = "9a71ffcd-7f60-46ca-8938-172a2aca64aa.md"
$wsDe.Cells.Item(3, 2).Value = ".md"
$wsDe.Cells.Item(3, 3).Value = "Ready for handoff"
$wsDe.Cells.Item(3, 4).Value = "e2e"
$wsDe.Cells.Item(3, 5).Value = "ht"
$wsDe.Cells.Item(3, 6).Value = "'False"
$wsDe.Cells.Item(3, 7).Value = "9a71ffcd-7f60-46ca-8938-172a2aca64aa.842dd4dc5e93682265b9bcceb6d8905f68e54fe8.de-de.xlf"
$wsDe.Cells.Item(3, 8).Value = "2016-08-13 14:48:23"
$wsDe.Cells.Item(3, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(3, 11).Value = "0001-01-01 00:00:00"
$wsDe.Cells.Item(3, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(3, 13).Value = "'True"
$wsDe.Cells.Item(3, 15).Value = "'False"

$wsDe.Hyperlinks.Add($wsDe.Cells.Item(3, 1), "https://github.com/OpenLocalizationTestOrg/oltest/blob/24dad09822a2e81cfe6c5f481aba677f2333a927/e2e/9a71ffcd-7f60-46ca-8938-172a2aca64aa.md", "", "", "9a71ffcd-7f60-46ca-8938-172a2aca64aa.md")
